# Updates "想去人数" (want-to-go count) values in column F across all sheets
# to reflect newly generated data (output generated at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1361
$ws.Range("F6").Value = 2638
$ws.Range("F7").Value = 1003
$ws.Range("F8").Value = 19111
$ws.Range("F9").Value = 71
$ws.Range("F10").Value = 2113
$ws.Range("F11").Value = 712
$ws.Range("F12").Value = 611
$ws.Range("F13").Value = 388
$ws.Range("F14").Value = 652
$ws.Range("F15").Value = 214
$ws.Range("F16").Value = 227
$ws.Range("F18").Value = 337
$ws.Range("F20").Value = 230
$ws.Range("F22").Value = 148

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 181
$ws.Range("F7").Value = 7
$ws.Range("F8").Value = 116
$ws.Range("F9").Value = 250
$ws.Range("F10").Value = 250
$ws.Range("F18").Value = 45
$ws.Range("F19").Value = 3

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 609

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 609
$ws.Range("F8").Value = 1361
$ws.Range("F11").Value = 181
$ws.Range("F13").Value = 2638
$ws.Range("F14").Value = 1003
$ws.Range("F15").Value = 19111
$ws.Range("F17").Value = 7
$ws.Range("F18").Value = 71
$ws.Range("F19").Value = 116
$ws.Range("F20").Value = 250
$ws.Range("F21").Value = 250
$ws.Range("F22").Value = 2113
$ws.Range("F23").Value = 713
$ws.Range("F25").Value = 388
$ws.Range("F26").Value = 653
$ws.Range("F27").Value = 214
$ws.Range("F28").Value = 227
$ws.Range("F33").Value = 337
$ws.Range("F36").Value = 230
$ws.Range("F39").Value = 148
$ws.Range("F40").Value = 45
$ws.Range("F41").Value = 3
